$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.637.12'
$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("D3").Value = '2.703.97'
$ws.Range("E3").Value = '  +2.26%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.93%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.545'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").Value = '2.702.96'
$ws.Range("E9").Value = '  +2.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.141'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.35%  '

$ws.Range("E11").Value = '  -0.31%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.31'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.16%  '

$ws.Range("E13").Value = '  +2.57%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.34'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.24%  '

$ws.Range("D15").Value = '3.194.87'
$ws.Range("E15").Value = '  +2.23%  '

$ws.Range("E16").Value = '  -0.25%  '

$ws.Range("D17").Value = '68.600.13'
$ws.Range("E17").Value = '  +0.71%  '

$ws.Range("D18").Value = '2.713.70'
$ws.Range("E18").Value = '  +2.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.37%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '366.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.64'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.83%  '

$ws.Range("E22").Value = '  +2.89%  '

$ws.Range("E23").Value = '  +2.43%  '

$ws.Range("E24").Value = '  +2.95%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.78%  '

$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.66%  '

$ws.Range("D28").Value = '2.827.38'
$ws.Range("E28").Value = '  +1.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000106'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.45%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.35%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '578.63'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.39%  '

$ws.Range("E32").Value = '  +2.76%  '

$ws.Range("E33").Value = '  +3.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.96'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.16%  '

$ws.Range("E35").Value = '  +3.78%  '

$ws.Range("E36").Value = '  +6.49%  '

$ws.Range("E37").Value = '  +0.00%  '

$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.90'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.81%  '

$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '160.70'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.37%  '

$ws.Range("E40").Value = '  +2.32%  '

$ws.Range("E41").Value = '  +2.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.74%  '

$ws.Range("E43").Value = '  +3.77%  '

$ws.Range("E44").Value = '  +0.36%  '

$ws.Range("E45").Value = '  -4.75%  '

$ws.Range("E46").Value = '  +0.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '157.95'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.68%  '

$ws.Range("E48").Value = '  +6.55%  '

$ws.Range("E49").Value = '  +5.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.601'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.90%  '

$ws.Range("E51").Value = '  -0.33%  '
